$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two fresh rows right after the current "Docentes responsáveis:"
#    row (row 12), before the current "Programa resumido:" row (row 13).
#    This makes room for two new rows that will hold the two docente names
#    (previously mis-placed a couple of rows further down), and shifts the
#    rest of the form (Programa resumido / Programa / Avaliação block /
#    Bibliografia / Requisitos) down by two rows, matching the target
#    layout. Inserting at row 13 twice pushes everything at/after row 13
#    down by two rows total.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# ---------------------------------------------------------------------------
# Helper pattern used below for every new/changed text cell in columns B/C:
#   - copy number/alignment/font formatting from a same-column cell that
#     already carries the right style (so no new style entries are minted)
#   - then write the literal text value
#   - finally, for any row whose column-A cell got auto-materialised by the
#     row insert (inherits the bold "label" style from column A) but should
#     stay empty, clear it explicitly.
# ---------------------------------------------------------------------------

# --- Row 10: "Objetivos:" value (column B/C) — replace placeholder text ----
$ws.Cells.Item(10,2).Value = "Adquirir conhecimentos dos principais processos de soldagem, nomenclatura das juntas soldadas, metalurgia física da soldagem, aspectos relativos à segurança e aplicações da soldagem em engenharia.Processos de metalurgia do pó metálico, de seus principais aspectos metalúrgicos, propriedades, aplicações, vantagens e desvantagens técnicas e econômicas. Identificação dos problemas comuns em componentes metálicos fundidos, soldados e sinterizados. Introdução à Manufatura Aditiva: Potencialidade e Técnicas."
$ws.Cells.Item(10,3).Value = "Adquirir conhecimentos dos principais processos de soldagem, nomenclatura das juntas soldadas, metalurgia física da soldagem, aspectos relativos à segurança e aplicações da soldagem em engenharia.Processos de metalurgia do pó metálico, de seus principais aspectos metalúrgicos, propriedades, aplicações, vantagens e desvantagens técnicas e econômicas. Identificação dos problemas comuns em componentes metálicos fundidos, soldados e sinterizados. Introdução à Manufatura Aditiva: Potencialidade e Técnicas."

# --- New row 13: docente #1 (no column-A label) ----------------------------
$ws.Cells.Item(8,2).Copy()
$ws.Cells.Item(13,2).PasteSpecial(-4122)
$ws.Cells.Item(8,3).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4122)
$ws.Cells.Item(13,2).Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Cells.Item(13,3).Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Cells.Item(13,1).Clear()

# --- New row 14: docente #2 (no column-A label) ----------------------------
$ws.Cells.Item(8,2).Copy()
$ws.Cells.Item(14,2).PasteSpecial(-4122)
$ws.Cells.Item(8,3).Copy()
$ws.Cells.Item(14,3).PasteSpecial(-4122)
$ws.Cells.Item(14,2).Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Cells.Item(14,3).Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Cells.Item(14,1).Clear()

# --- Row 15 (shifted "Programa resumido:"): correct summary text ----------
$ws.Cells.Item(15,2).Value = "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS; 2. PROCESSOS DE SOLDAGEM; 3. NOMENCLATURA DAS JUNTAS SOLDADAS; 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS; 5. SEGURANÇA NO PROCESSO DE SOLDAGEM; 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA; 7.PÓS METALICOS – OBTENÇÃO, CARACTERIZAÇÃO E APLICAÇÃO NA METALURGIA DO PÓ. 8. TÉCNICAS DE MISTURA, 9. PROCESSOS DE FABRICAÇÃO DE PEÇAS VERDES, 10. SINTERIZAÇÃO, 11. UTILIZAÇÃO DO LASER E DE FEIXE DE ELÉTRONS12. PRÁTICA EXPERIMENTAL SUPERVISIONADA."
$ws.Cells.Item(15,3).Value = "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS; 2. PROCESSOS DE SOLDAGEM; 3. NOMENCLATURA DAS JUNTAS SOLDADAS; 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS; 5. SEGURANÇA NO PROCESSO DE SOLDAGEM; 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA; 7.PÓS METALICOS – OBTENÇÃO, CARACTERIZAÇÃO E APLICAÇÃO NA METALURGIA DO PÓ. 8. TÉCNICAS DE MISTURA, 9. PROCESSOS DE FABRICAÇÃO DE PEÇAS VERDES, 10. SINTERIZAÇÃO, 11. UTILIZAÇÃO DO LASER E DE FEIXE DE ELÉTRONS12. PRÁTICA EXPERIMENTAL SUPERVISIONADA."

# --- Row 17 (shifted "Programa:"): correct full program text --------------
$ws.Cells.Item(17,2).Value = "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS: Razões técnicas para a junção de materiais, junção por difusão, brasagem, soldagem por explosão, elementos de fixação. 2. PROCESSOS DE SOLDAGEM: Definição de soldagem por fusão, física da soldagem, principais processos de soldagem. 3. NOMENCLATURA DAS JUNTAS SOLDADAS: Desenho e simbologia para soldagem, símbolos básicos, tipos de juntas e soldas, simbologia para soldas em desenho. 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS: metalurgia da soldagem, estruturas de solidificação, transformações de fase pós-soldagem, transformações de fases em juntas de aço soldadas, ligas de alumínio, ligas de cobre e em metais e ligas especiais. 5. SEGURANÇA NO PROCESSO DE SOLDAGEM: Problemas associados à vaporização de metais, luminosidade, calor e eletricidade. 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA: Exemplos de estruturas soldadas em engenharia, descontinuidades e defeitos de soldagem, métodos de inspeção em soldas. 7. PRÁTICA EXPERIMENTAL SUPERVISIONADA: Caracterização microestrutural de juntas soldadas (materiais e processos a serem definidos na ocasião da prática experimental), incluindo a redação de relatório técnico de cada grupo. 8. Pós Metálicos - obtenção por processos químicos, termoquímicos, eletrolíticos, atomização e moagem, Caracterização de pós e sua aplicação na metalurgia do pó.9. Técnicas de mistura, aspectos sobre o transporte e armazenamento de pós, 10-Processos de fabricação de peças `"verdes`" por compactação uniaxial e isostática, 11- Técnicas de sinterização e fenômenos envolvidos, 12 Sinterização/refusão a LASER para prototipagem rápida (impressão 3D). Feixe de elétrons: obtenção e aplicações."
$ws.Cells.Item(17,3).Value = "1. TÉCNICAS DE JUNÇÃO DE MATERIAIS: Razões técnicas para a junção de materiais, junção por difusão, brasagem, soldagem por explosão, elementos de fixação. 2. PROCESSOS DE SOLDAGEM: Definição de soldagem por fusão, física da soldagem, principais processos de soldagem. 3. NOMENCLATURA DAS JUNTAS SOLDADAS: Desenho e simbologia para soldagem, símbolos básicos, tipos de juntas e soldas, simbologia para soldas em desenho. 4. METALURGIA FÍSICA DAS REGIÕES SOLDADAS: metalurgia da soldagem, estruturas de solidificação, transformações de fase pós-soldagem, transformações de fases em juntas de aço soldadas, ligas de alumínio, ligas de cobre e em metais e ligas especiais. 5. SEGURANÇA NO PROCESSO DE SOLDAGEM: Problemas associados à vaporização de metais, luminosidade, calor e eletricidade. 6. APLICAÇÕES DE JUNTAS SOLDADAS EM ENGENHARIA: Exemplos de estruturas soldadas em engenharia, descontinuidades e defeitos de soldagem, métodos de inspeção em soldas. 7. PRÁTICA EXPERIMENTAL SUPERVISIONADA: Caracterização microestrutural de juntas soldadas (materiais e processos a serem definidos na ocasião da prática experimental), incluindo a redação de relatório técnico de cada grupo. 8. Pós Metálicos - obtenção por processos químicos, termoquímicos, eletrolíticos, atomização e moagem, Caracterização de pós e sua aplicação na metalurgia do pó.9. Técnicas de mistura, aspectos sobre o transporte e armazenamento de pós, 10-Processos de fabricação de peças `"verdes`" por compactação uniaxial e isostática, 11- Técnicas de sinterização e fenômenos envolvidos, 12 Sinterização/refusão a LASER para prototipagem rápida (impressão 3D). Feixe de elétrons: obtenção e aplicações."

# --- Row 20 (shifted "Método:"): correct evaluation method text -----------
$ws.Cells.Item(20,2).Value = "O aluno será avaliado por duas avaliações, sendo que a segunda avaliação terá peso 2."
$ws.Cells.Item(20,3).Value = "O aluno será avaliado por duas avaliações, sendo que a segunda avaliação terá peso 2."

# --- Row 21 (shifted "Critério:"): correct grading formula text -----------
$ws.Cells.Item(21,2).Value = "Nota Final NF = [Avaliação 1 + 2*(Avaliação 2)]/3"
$ws.Cells.Item(21,3).Value = "Nota Final NF = [Avaliação 1 + 2*(Avaliação 2)]/3"

# --- Row 22 (shifted "Norma de recuperação:"): correct makeup-exam text ---
$ws.Cells.Item(22,2).Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Cells.Item(22,3).Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."

# --- Row 23 (shifted "Bibliografia:"): correct bibliography text ----------
$ws.Cells.Item(23,2).Value = "1. WAINER, E. et al. Soldagem - Processos e Metalurgia. São Paulo: Edgar Blücher, 1992. 494 p. 2. QUITES, A. M., DUTRA, J. C. Tecnologia da soldagem a arco voltaico. Florianópolis: EDEME, 1979. 248 p. 3. GOURD, L. M. Principles of welding technology. London: Edward Arnold, 1980. 218 p. 4. KOU, S. Welding metallurgy, 2nd ed.: John Wiley & Sons, 2003. 461 p. 5. MESSLER, Jr. R. W. Principles of welding: Processes, physics, chemistry and metallurgy: Wiley VCH Verlag GmbH & Co., 2004. 662 p.6. KALPAKJIAN, S.; SCHMID, S. Manufacturing processes for engineering materials. 5ª ed., Pearson Education, New Jersey, 2007.7. GERMAN, R.M. Sintering theory and practice. New York, Wiley-Interscience, 19968. GIBSON, I., ROSEN, D., STUCKER, B., Additive Manufacturing Technologies, New York, Springer Verlag, 2015."
$ws.Cells.Item(23,3).Value = "1. WAINER, E. et al. Soldagem - Processos e Metalurgia. São Paulo: Edgar Blücher, 1992. 494 p. 2. QUITES, A. M., DUTRA, J. C. Tecnologia da soldagem a arco voltaico. Florianópolis: EDEME, 1979. 248 p. 3. GOURD, L. M. Principles of welding technology. London: Edward Arnold, 1980. 218 p. 4. KOU, S. Welding metallurgy, 2nd ed.: John Wiley & Sons, 2003. 461 p. 5. MESSLER, Jr. R. W. Principles of welding: Processes, physics, chemistry and metallurgy: Wiley VCH Verlag GmbH & Co., 2004. 662 p.6. KALPAKJIAN, S.; SCHMID, S. Manufacturing processes for engineering materials. 5ª ed., Pearson Education, New Jersey, 2007.7. GERMAN, R.M. Sintering theory and practice. New York, Wiley-Interscience, 19968. GIBSON, I., ROSEN, D., STUCKER, B., Additive Manufacturing Technologies, New York, Springer Verlag, 2015."

Write-Output "done"
